$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.954.39"
$ws.Range("E2").Value = "  +2.86%  "
$ws.Range("D3").Value = "2.637.27"
$ws.Range("E3").Value = "  +9.61%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.94"
$ws.Range("E5").Value = "  +4.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.95"
$ws.Range("E6").Value = "  +6.57%  "
$ws.Range("E7").Value = "  +7.45%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  +15.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.44"
$ws.Range("E10").Value = "  +12.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "55.09"
$ws.Range("E11").Value = "  +2.27%  "
$ws.Range("E12").Value = "  +6.78%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.39"
$ws.Range("E13").Value = "  +16.84%  "
$ws.Range("D14").Value = "3.040.53"
$ws.Range("E14").Value = "  +10.28%  "
$ws.Range("E15").Value = "  +2.04%  "
$ws.Range("D16").Value = "2.648.82"
$ws.Range("E16").Value = "  +9.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.933"
$ws.Range("E17").Value = "  +9.94%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "15.24"
$ws.Range("E18").Value = "  +6.50%  "
$ws.Range("D19").Value = "47.429.61"
$ws.Range("E19").Value = "  +4.08%  "
$ws.Range("E20").Value = "  +8.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.39"
$ws.Range("E21").Value = "  +3.10%  "
$ws.Range("E22").Value = "  +9.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.56"
$ws.Range("E23").Value = "  +6.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "259.93"
$ws.Range("E24").Value = "  +6.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.14"
$ws.Range("E25").Value = "  +11.00%  "
$ws.Range("E26").Value = "  +16.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "30.46"
$ws.Range("E27").Value = "  +43.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "41.72"
$ws.Range("E29").Value = "  +8.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.74"
$ws.Range("E30").Value = "  +9.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.31"
$ws.Range("E31").Value = "  +4.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.30"
$ws.Range("E32").Value = "  +13.75%  "
$ws.Range("E33").Value = "  -1.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.32"
$ws.Range("E34").Value = "  +16.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.89"
$ws.Range("E35").Value = "  +5.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0847"
$ws.Range("E36").Value = "  +8.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "153.46"
$ws.Range("E37").Value = "  +2.90%  "
$ws.Range("E38").Value = "  +4.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.124"
$ws.Range("E39").Value = "  +6.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.98"
$ws.Range("E40").Value = "  +10.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.38"
$ws.Range("E41").Value = "  +12.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.74"
$ws.Range("E42").Value = "  +14.35%  "
$ws.Range("E43").Value = "  +9.87%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.53"
$ws.Range("E44").Value = "  +38.13%  "
$ws.Range("D45").Value = "2.055.77"
$ws.Range("E45").Value = "  +5.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.998"
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "93.33"
$ws.Range("E47").Value = "  +1.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "114.64"
$ws.Range("E48").Value = "  +11.30%  "
$ws.Range("E49").Value = "  +4.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.30"
$ws.Range("E50").Value = "  +6.58%  "
$ws.Range("E51").Value = "  +7.28%  "
